$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text so numeric-looking strings (e.g. "1.001") are not
# auto-converted to floating point numbers by Excel type inference.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '29.987.08'
$ws.Range('E2').Value = '  -1.03%  '
$ws.Range('D3').Value = '1.900.58'
$ws.Range('E3').Value = '  -1.11%  '
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  -0.18%  '
$ws.Range('D5').Value = '0.7419'
$ws.Range('E5').Value = '  +0.43%  '
$ws.Range('D6').Value = '241.60'
$ws.Range('E6').Value = '  -0.76%  '
$ws.Range('E7').Value = '  -0.24%  '
$ws.Range('D8').Value = '0.3063'
$ws.Range('E8').Value = '  -2.59%  '
$ws.Range('D9').Value = '25.38'
$ws.Range('E9').Value = '  -7.38%  '
$ws.Range('D10').Value = '0.06877'
$ws.Range('E10').Value = '  -1.71%  '
$ws.Range('D11').Value = '0.08001'
$ws.Range('E11').Value = '  -0.27%  '
$ws.Range('D12').Value = '0.7515'
$ws.Range('E12').Value = '  -2.52%  '
$ws.Range('D13').Value = '1.901.04'
$ws.Range('E13').Value = '  -0.09%  '
$ws.Range('D14').Value = '5.255'
$ws.Range('E14').Value = '  -1.47%  '
$ws.Range('D15').Value = '91.20'
$ws.Range('E15').Value = '  -1.75%  '
$ws.Range('D16').Value = '6.133'
$ws.Range('E16').Value = '  +5.01%  '
$ws.Range('D17').Value = '30.000.45'
$ws.Range('E17').Value = '  -0.78%  '
$ws.Range('D18').Value = '13.94'
$ws.Range('E18').Value = '  -2.81%  '
$ws.Range('D19').Value = '0.000007735'
$ws.Range('E19').Value = '  -1.97%  '
$ws.Range('D20').Value = '236.63'
$ws.Range('E20').Value = '  -5.20%  '
$ws.Range('B21').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C21').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D21').Value = '2.164.90'
$ws.Range('E21').Value = '  +2.38%  '
$ws.Range('B22').Value = 'Dai'
$ws.Range('C22').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D22').Value = '1.001'
$ws.Range('E22').Value = '  -0.26%  '
$ws.Range('E23').Value = '  -0.15%  '
$ws.Range('D24').Value = '7.056'
$ws.Range('E24').Value = '  +6.78%  '
$ws.Range('B25').Value = 'Monero'
$ws.Range('C25').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D25').Value = '167.48'
$ws.Range('E25').Value = '  +1.27%  '
$ws.Range('B26').Value = 'Cosmos'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D26').Value = '9.302'
$ws.Range('E26').Value = '  -2.10%  '
$ws.Range('D27').Value = '18.74'
$ws.Range('E27').Value = '  -1.00%  '
$ws.Range('D28').Value = '0.1260'
$ws.Range('E28').Value = '  -1.76%  '
$ws.Range('D29').Value = '2.041'
$ws.Range('E29').Value = '  -5.43%  '
$ws.Range('D30').Value = '1.350'
$ws.Range('E30').Value = '  -0.54%  '
$ws.Range('D31').Value = '1.525'
$ws.Range('E31').Value = '  -2.37%  '
$ws.Range('D32').Value = '4.296'
$ws.Range('E32').Value = '  -2.18%  '
$ws.Range('D33').Value = '4.035'
$ws.Range('E33').Value = '  -1.21%  '
$ws.Range('D34').Value = '0.05294'
$ws.Range('E34').Value = '  +1.69%  '
$ws.Range('D35').Value = '1.277'
$ws.Range('E35').Value = '  -1.38%  '
$ws.Range('D36').Value = '0.7356'
$ws.Range('E36').Value = '  -2.08%  '
$ws.Range('D37').Value = '2.724'
$ws.Range('E37').Value = '  -1.62%  '
$ws.Range('D38').Value = '0.01941'
$ws.Range('E38').Value = '  +0.17%  '
$ws.Range('D39').Value = '2.773'
$ws.Range('E39').Value = '  -0.30%  '
$ws.Range('D40').Value = '6.228'
$ws.Range('E40').Value = '  -4.00%  '
$ws.Range('D41').Value = '0.4441'
$ws.Range('E41').Value = '  -0.86%  '
$ws.Range('D42').Value = '72.54'
$ws.Range('E42').Value = '  -4.74%  '
$ws.Range('D43').Value = '1.936'
$ws.Range('E43').Value = '  -0.05%  '
$ws.Range('E44').Value = '  -0.17%  '
$ws.Range('B45').Value = 'TrustWalletToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D45').Value = '0.8321'
$ws.Range('E45').Value = '  -0.78%  '
$ws.Range('B46').Value = 'Aptos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D46').Value = '7.735'
$ws.Range('E46').Value = '  +1.44%  '
$ws.Range('D47').Value = '101.44'
$ws.Range('E47').Value = '  +0.29%  '
$ws.Range('D48').Value = '9.860'
$ws.Range('E48').Value = '  +0.08%  '
$ws.Range('D49').Value = '2.062.39'
$ws.Range('E49').Value = '  +0.61%  '
$ws.Range('D50').Value = '36.53'
$ws.Range('E50').Value = '  -2.12%  '
$ws.Range('D51').Value = '0.05977'
$ws.Range('E51').Value = '  -0.50%  '

# Restore default (General) formatting on column D now that the text values
# are stored, so no stray style indices remain on the cells.
$ws.Range("D2:D51").ClearFormats()
